$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "Licia Fondriest "
$ws.Range("B20").Value = "Stefano Tita | Clitoriders"
$ws.Range("C20").Value = "Stefano  Galvagni | Clitoriders"
$ws.Range("D20").Value = "Edoardo Pomarolli | Modium"
$ws.Range("E20").Value = "Luca Giordani | Shark Attack"
$ws.Range("F20").Value = "Halzyd  Pupuleku | F.C. Sala Giardini"
